$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.701.81"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "3.558.00"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "574.21"
$ws.Range("E5").Value = "  -3.50%  "
$ws.Range("D6").Value = "187.66"
$ws.Range("E6").Value = "  -2.83%  "
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  -3.22%  "
$ws.Range("D8").Value = "3.553.32"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  -3.92%  "
$ws.Range("D11").Value = "0.656"
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("D12").Value = "55.78"
$ws.Range("E12").Value = "  -4.37%  "
$ws.Range("D13").Value = "0.0000298"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").Value = "9.66"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").Value = "4.131.74"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "19.80"
$ws.Range("E16").Value = "  +2.08%  "
$ws.Range("D17").Value = "3.560.62"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").Value = "69.583.15"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "12.47"
$ws.Range("E19").Value = "  -1.95%  "
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").Value = "1.03"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("D22").Value = "470.35"
$ws.Range("E22").Value = "  -6.26%  "
$ws.Range("D23").Value = "19.21"
$ws.Range("E23").Value = "  +12.99%  "
$ws.Range("D24").Value = "5.04"
$ws.Range("E24").Value = "  -8.06%  "
$ws.Range("D25").Value = "4.33"
$ws.Range("E25").Value = "  -3.34%  "
$ws.Range("D26").Value = "88.11"
$ws.Range("E26").Value = "  -3.54%  "
$ws.Range("D27").Value = "3.03"
$ws.Range("E27").Value = "  -2.00%  "
$ws.Range("D28").Value = "10.91"
$ws.Range("E28").Value = "  -2.86%  "
$ws.Range("D29").Value = "9.32"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("D30").Value = "31.95"
$ws.Range("E30").Value = "  -1.51%  "
$ws.Range("D31").Value = "7.61"
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("D33").Value = "12.03"
$ws.Range("E33").Value = "  -1.46%  "
$ws.Range("D34").Value = "65.53"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "581.24"
$ws.Range("E35").Value = "  -5.67%  "
$ws.Range("D36").Value = "38.46"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "0.0₃0795"
$ws.Range("E38").Value = "  -5.33%  "
$ws.Range("D39").Value = "0.394"
$ws.Range("E39").Value = "  -1.71%  "
$ws.Range("D40").Value = "0.139"
$ws.Range("E40").Value = "  -6.11%  "
$ws.Range("D41").Value = "3.48"
$ws.Range("E41").Value = "  -5.16%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.212.18"
$ws.Range("E42").Value = "  -4.27%  "
$ws.Range("D43").Value = "2.84"
$ws.Range("E43").Value = "  +5.02%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "3.14"
$ws.Range("E44").Value = "  +10.99%  "
$ws.Range("D45").Value = "3.10"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").Value = "0.0441"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("D47").Value = "9.39"
$ws.Range("E47").Value = "  +2.95%  "
$ws.Range("E48").Value = "  +1.12%  "
$ws.Range("D49").Value = "0.136"
$ws.Range("E49").Value = "  -1.08%  "
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("D51").Value = "137.32"
$ws.Range("E51").Value = "  -3.17%  "

Write-Host "Applied $(97) changes"
